# Center the "Figure" paragraph style, per the updated NAFO plotting
# guidelines (figures should be centered on the page).
$d = $word.ActiveDocument

$figureStyle = $d.Styles("Figure")
$figureStyle.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter
